$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = "Normal"
}

Set-TextValue 'D2' '60.994.75'
$ws.Range('E2').Value = '  +0.18%  '
Set-TextValue 'D3' '2.924.82'
$ws.Range('E3').Value = '  +0.20%  '
$ws.Range('E4').Value = '  -0.03%  '
Set-TextValue 'D5' '590.88'
$ws.Range('E5').Value = '  +1.00%  '
Set-TextValue 'D6' '147.12'
$ws.Range('E6').Value = '  +1.12%  '
$ws.Range('E7').Value = '  +0.00%  '
$ws.Range('E8').Value = '  +0.54%  '
$ws.Range('E9').Value = '  +0.53%  '
$ws.Range('E10').Value = '  -0.43%  '
$ws.Range('E11').Value = '  -1.20%  '
$ws.Range('E12').Value = '  +0.06%  '
Set-TextValue 'D13' '33.71'
$ws.Range('E13').Value = '  +0.13%  '
Set-TextValue 'D15' '3.409.34'
$ws.Range('E15').Value = '  +0.19%  '
Set-TextValue 'D16' '60.928.29'
$ws.Range('E16').Value = '  +0.11%  '
$ws.Range('B17').Value = 'Polkadot'
$ws.Range('C17').Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue 'D17' '6.71'
$ws.Range('E17').Value = '  -0.70%  '
$ws.Range('B18').Value = 'WrappedEther'
$ws.Range('C18').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue 'D18' '2.924.75'
$ws.Range('E18').Value = '  +0.17%  '
Set-TextValue 'D19' '432.18'
$ws.Range('E19').Value = '  +0.24%  '
$ws.Range('E20').Value = '  -1.56%  '
Set-TextValue 'D21' '0.680'
$ws.Range('E21').Value = '  -0.44%  '
$ws.Range('E22').Value = '  -0.66%  '
Set-TextValue 'D23' '81.43'
$ws.Range('E23').Value = '  +1.29%  '
$ws.Range('E24').Value = '  +0.73%  '
$ws.Range('E25').Value = '  -0.41%  '
Set-TextValue 'D26' '11.90'
$ws.Range('E26').Value = '  -0.33%  '
$ws.Range('E27').Value = '  -0.05%  '
$ws.Range('E28').Value = '  +4.87%  '
$ws.Range('E29').Value = '  -0.03%  '
$ws.Range('E30').Value = '  -3.11%  '
$ws.Range('E31').Value = '  +0.43%  '
$ws.Range('E32').Value = '  +2.25%  '
$ws.Range('E33').Value = '  -0.02%  '
$ws.Range('E34').Value = '  -0.86%  '
$ws.Range('E35').Value = '  +0.14%  '
$ws.Range('E36').Value = '  -0.11%  '
Set-TextValue 'D37' '3.01'
$ws.Range('E37').Value = '  -1.16%  '
$ws.Range('E38').Value = '  -1.16%  '
$ws.Range('E39').Value = '  -5.21%  '
$ws.Range('E40').Value = '  -1.02%  '
Set-TextValue 'D41' '41.43'
$ws.Range('E41').Value = '  +1.45%  '
$ws.Range('E42').Value = '  -4.49%  '
Set-TextValue 'D43' '377.99'
$ws.Range('E43').Value = '  +0.04%  '
Set-TextValue 'D44' '2.707.28'
$ws.Range('E44').Value = '  +1.00%  '
$ws.Range('E45').Value = '  -1.65%  '
Set-TextValue 'D46' '134.15'
$ws.Range('E46').Value = '  +1.11%  '
Set-TextValue 'D48' '23.84'
$ws.Range('E48').Value = '  -4.24%  '
$ws.Range('E49').Value = '  -0.62%  '
$ws.Range('E50').Value = '  -2.77%  '
$ws.Range('E51').Value = '  -0.68%  '
